# DaySale report update: add two new low-stock items (GAVISCON, PANADOL),
# bump the totals/footer accordingly and refresh the "printed at" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room: insert two new data rows right above the current row 9
#    (which holds "SPASMOFEN..."). This pushes the existing row 9/10
#    (SPASMOFEN / syringes) down to rows 11/12, and the footer rows
#    (11/12) down to rows 13/14 - exactly like the target workbook.
# ---------------------------------------------------------------------
$ws.Rows("9:10").Insert()

# Copy the formatting (styles/number-formats/borders/merges-pattern) of
# row 8 onto the two freshly inserted blank rows, so they look identical
# to the other product rows.
$ws.Range("A8:Q8").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4122)
$ws.Range("A9:Q9").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-create the cell merges for the two new rows (same pattern used by
# every other data row in the table).
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

# Match the row heights of the final report (row 9 / row 10 keep the
# same heights the two rows already had in this slot; row 12, holding
# the item that got pushed down, ends up re-measured at 25.5 as well).
$ws.Rows("9").RowHeight = 25.5
$ws.Rows("10").RowHeight = 24.75
$ws.Rows("12").RowHeight = 25.5

# ---------------------------------------------------------------------
# 2) Fill in the new row 9 - GAVISCON LIQUID 24 SACHETS 10 ML
# ---------------------------------------------------------------------
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "GAVISCON LIQUID 24 SACHETS 10 ML"
$ws.Range("H9").Value = "0:12"
$ws.Range("L9").Value = "'1"
$ws.Range("N9").Value = "'288.00"
$ws.Range("P9").Value = "'11.5200"
$ws.Range("Q9").Value = "0:1"

# ---------------------------------------------------------------------
# 3) Fill in the new row 10 - PANADOL ADVANCE 500 MG 48 TABLETS
# ---------------------------------------------------------------------
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "PANADOL ADVANCE 500 MG 48 TABLETS"
$ws.Range("H10").Value = "2:3"
$ws.Range("L10").Value = "'1"
$ws.Range("N10").Value = "'92.00"
$ws.Range("P10").Value = "'23.0000"
$ws.Range("Q10").Value = "0:1"

# Re-apply the reference formatting on the text cells we just overwrote
# with a leading quote (forces text storage) so their number format /
# style stays identical to the sibling rows.
$ws.Range("L8:Q8").Copy()
$ws.Range("L9:Q9").PasteSpecial(-4122)
$ws.Range("L9:Q9").Copy()
$ws.Range("L10:Q10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Values above got blanked by the formats-only paste? No - PasteSpecial
# with "paste formats" does not touch the values, only styles, so the
# text entered above survives untouched.

# ---------------------------------------------------------------------
# 4) Renumber the rows that shifted down (old row 9 -> 11, old row 10 -> 12)
# ---------------------------------------------------------------------
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6

# ---------------------------------------------------------------------
# 5) Update the totals footer (old row 11 -> now row 13) and the
#    timestamp / footer line (old row 12 -> now row 14).
# ---------------------------------------------------------------------
$ws.Range("P13").Value = 110.19
$ws.Range("A14").Value = "Friday, 22 August, 2025 5:19 PM"

Write-Host "Edit complete"
